$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph.
$verFind = $d.Content
$verFound = $verFind.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verIndex = $verFind.Paragraphs.First.Index

# Locate the "© 2020 . Contact: ..." paragraph right after it.
$copyFind = $d.Content
$copyFound = $copyFind.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyIndex = $copyFind.Paragraphs.First.Index

# The empty paragraph immediately preceding "Ver no Jupiter ..." is removed too,
# leaving the other blank paragraph (the one right before the page break) intact.
$emptyPara = $d.Paragraphs.Item($verIndex - 1)
$copyPara = $d.Paragraphs.Item($copyIndex)

$delRange = $d.Range($emptyPara.Range.Start, $copyPara.Range.End)
$delRange.Delete()
